$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Make room for the new row (Bid ID 10 gets a second split, "Facility5"/B) by
# inserting a fresh row 15 before the sheet's former last row (14) grows.
$ws.Rows.Item(15).Insert()

# Row 2
$ws.Cells.Item(2, 6).Value = 35000
$cell = $ws.Cells.Item(2, 7)
$cell.NumberFormat = "@"
$cell.Value = "A"
$ws.Cells.Item(2, 8).Value = 20
$cell = $ws.Cells.Item(2, 9)
$cell.NumberFormat = "@"
$cell.Value = "1%"
$ws.Cells.Item(2, 10).Value = 19.8
$ws.Cells.Item(2, 11).Value = 6930
$ws.Cells.Item(2, 12).Value = 350
$ws.Cells.Item(2, 13).Value = 28070
$cell = $ws.Cells.Item(2, 14)
$cell.NumberFormat = "@"
$cell.Value = "0%"
$ws.Cells.Item(2, 15).Value = 0

# Row 3
$ws.Cells.Item(3, 6).Value = 35000
$cell = $ws.Cells.Item(3, 7)
$cell.NumberFormat = "@"
$cell.Value = "B"
$ws.Cells.Item(3, 8).Value = 60
$cell = $ws.Cells.Item(3, 9)
$cell.NumberFormat = "@"
$cell.Value = "3%"
$ws.Cells.Item(3, 10).Value = 58.2
$ws.Cells.Item(3, 11).Value = 20370
$ws.Cells.Item(3, 12).Value = 350
$ws.Cells.Item(3, 13).Value = 14630
$cell = $ws.Cells.Item(3, 14)
$cell.NumberFormat = "@"
$cell.Value = "5%"
$ws.Cells.Item(3, 15).Value = 1018.5

# Row 5
$ws.Cells.Item(5, 6).Value = 583752
$ws.Cells.Item(5, 11).Value = 259320.6
$ws.Cells.Item(5, 12).Value = 3742
$ws.Cells.Item(5, 13).Value = 324431.4

# Row 6
$ws.Cells.Item(6, 1).Value = 2
$cell = $ws.Cells.Item(6, 2)
$cell.NumberFormat = "@"
$cell.Value = "C"
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "B"
$ws.Cells.Item(6, 5).Value = 156
$ws.Cells.Item(6, 6).Value = 40248
$cell = $ws.Cells.Item(6, 7)
$cell.NumberFormat = "@"
$cell.Value = "B"
$ws.Cells.Item(6, 8).Value = 80
$cell = $ws.Cells.Item(6, 9)
$cell.NumberFormat = "@"
$cell.Value = "3%"
$ws.Cells.Item(6, 10).Value = 77.59999999999999
$ws.Cells.Item(6, 11).Value = 20020.8
$ws.Cells.Item(6, 12).Value = 258
$ws.Cells.Item(6, 13).Value = 20227.2
$cell = $ws.Cells.Item(6, 14)
$cell.NumberFormat = "@"
$cell.Value = "5%"
$ws.Cells.Item(6, 15).Value = 1001.04

# Row 7
$ws.Cells.Item(7, 1).Value = 3
$cell = $ws.Cells.Item(7, 3)
$cell.NumberFormat = "@"
$cell.Value = "Facility1"
$ws.Cells.Item(7, 5).Value = 423
$ws.Cells.Item(7, 6).Value = 253800
$ws.Cells.Item(7, 8).Value = 60
$ws.Cells.Item(7, 10).Value = 57.59999999999999
$ws.Cells.Item(7, 11).Value = 34560
$ws.Cells.Item(7, 12).Value = 600
$ws.Cells.Item(7, 13).Value = 219240
$ws.Cells.Item(7, 15).Value = 2419.2

# Row 8
$cell = $ws.Cells.Item(8, 2)
$cell.NumberFormat = "@"
$cell.Value = "A"
$ws.Cells.Item(8, 6).Value = 1359000
$cell = $ws.Cells.Item(8, 7)
$cell.NumberFormat = "@"
$cell.Value = "C"
$ws.Cells.Item(8, 8).Value = 19
$cell = $ws.Cells.Item(8, 9)
$cell.NumberFormat = "@"
$cell.Value = "4%"
$ws.Cells.Item(8, 10).Value = 18.24
$ws.Cells.Item(8, 11).Value = 54719.99999999999
$ws.Cells.Item(8, 12).Value = 3000
$ws.Cells.Item(8, 13).Value = 1304280
$cell = $ws.Cells.Item(8, 14)
$cell.NumberFormat = "@"
$cell.Value = "7%"
$ws.Cells.Item(8, 15).Value = 3830.4

# Row 9
$ws.Cells.Item(9, 1).Value = 4
$cell = $ws.Cells.Item(9, 2)
$cell.NumberFormat = "@"
$cell.Value = "B"
$ws.Cells.Item(9, 5).Value = 453
$ws.Cells.Item(9, 6).Value = 1209510
$cell = $ws.Cells.Item(9, 7)
$cell.NumberFormat = "@"
$cell.Value = "A"
$ws.Cells.Item(9, 8).Value = 23
$cell = $ws.Cells.Item(9, 9)
$cell.NumberFormat = "@"
$cell.Value = "1%"
$ws.Cells.Item(9, 10).Value = 22.77
$ws.Cells.Item(9, 11).Value = 60795.9
$ws.Cells.Item(9, 12).Value = 2670
$ws.Cells.Item(9, 13).Value = 1148714.1
$cell = $ws.Cells.Item(9, 14)
$cell.NumberFormat = "@"
$cell.Value = "0%"
$ws.Cells.Item(9, 15).Value = 0

# Row 10
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 5).Value = 342
$ws.Cells.Item(10, 6).Value = 15390
$cell = $ws.Cells.Item(10, 7)
$cell.NumberFormat = "@"
$cell.Value = "C"
$cell = $ws.Cells.Item(10, 9)
$cell.NumberFormat = "@"
$cell.Value = "4%"
$ws.Cells.Item(10, 10).Value = 23.04
$ws.Cells.Item(10, 11).Value = 1036.8
$ws.Cells.Item(10, 12).Value = 45
$ws.Cells.Item(10, 13).Value = 14353.2
$cell = $ws.Cells.Item(10, 14)
$cell.NumberFormat = "@"
$cell.Value = "7%"
$ws.Cells.Item(10, 15).Value = 72.57600000000001

# Row 11
$ws.Cells.Item(11, 1).Value = 6
$ws.Cells.Item(11, 5).Value = 653
$ws.Cells.Item(11, 6).Value = 158026
$cell = $ws.Cells.Item(11, 7)
$cell.NumberFormat = "@"
$cell.Value = "B"
$ws.Cells.Item(11, 8).Value = 24
$cell = $ws.Cells.Item(11, 9)
$cell.NumberFormat = "@"
$cell.Value = "3%"
$ws.Cells.Item(11, 10).Value = 23.28
$ws.Cells.Item(11, 11).Value = 5633.76
$ws.Cells.Item(11, 12).Value = 242
$ws.Cells.Item(11, 13).Value = 152392.24
$cell = $ws.Cells.Item(11, 14)
$cell.NumberFormat = "@"
$cell.Value = "5%"
$ws.Cells.Item(11, 15).Value = 281.688

# Row 12
$ws.Cells.Item(12, 1).Value = 7
$cell = $ws.Cells.Item(12, 3)
$cell.NumberFormat = "@"
$cell.Value = "Facility2"
$ws.Cells.Item(12, 5).Value = 432
$ws.Cells.Item(12, 6).Value = 286848
$cell = $ws.Cells.Item(12, 7)
$cell.NumberFormat = "@"
$cell.Value = "A"
$ws.Cells.Item(12, 8).Value = 23
$cell = $ws.Cells.Item(12, 9)
$cell.NumberFormat = "@"
$cell.Value = "1%"
$ws.Cells.Item(12, 10).Value = 22.77
$ws.Cells.Item(12, 11).Value = 15119.28
$ws.Cells.Item(12, 12).Value = 664
$ws.Cells.Item(12, 13).Value = 271728.72

# Row 13
$ws.Cells.Item(13, 1).Value = 8
$cell = $ws.Cells.Item(13, 3)
$cell.NumberFormat = "@"
$cell.Value = "Facility3"
$ws.Cells.Item(13, 5).Value = 456
$ws.Cells.Item(13, 6).Value = 10944
$cell = $ws.Cells.Item(13, 7)
$cell.NumberFormat = "@"
$cell.Value = "B"
$cell = $ws.Cells.Item(13, 9)
$cell.NumberFormat = "@"
$cell.Value = "3%"
$ws.Cells.Item(13, 10).Value = 12.61
$ws.Cells.Item(13, 11).Value = 302.64
$ws.Cells.Item(13, 12).Value = 24
$ws.Cells.Item(13, 13).Value = 10641.36
$cell = $ws.Cells.Item(13, 14)
$cell.NumberFormat = "@"
$cell.Value = "5%"
$ws.Cells.Item(13, 15).Value = 15.132

# Row 14
$ws.Cells.Item(14, 1).Value = 9
$cell = $ws.Cells.Item(14, 3)
$cell.NumberFormat = "@"
$cell.Value = "Facility4"
$ws.Cells.Item(14, 5).Value = 234
$ws.Cells.Item(14, 6).Value = 54288
$cell = $ws.Cells.Item(14, 7)
$cell.NumberFormat = "@"
$cell.Value = "C"
$cell = $ws.Cells.Item(14, 9)
$cell.NumberFormat = "@"
$cell.Value = "4%"
$ws.Cells.Item(14, 10).Value = 12.48
$ws.Cells.Item(14, 11).Value = 2895.36
$ws.Cells.Item(14, 12).Value = 232
$ws.Cells.Item(14, 13).Value = 51392.64
$cell = $ws.Cells.Item(14, 14)
$cell.NumberFormat = "@"
$cell.Value = "7%"
$ws.Cells.Item(14, 15).Value = 202.6752

# Row 15
$ws.Cells.Item(15, 1).Value = 10
$cell = $ws.Cells.Item(15, 2)
$cell.NumberFormat = "@"
$cell.Value = "A"
$cell = $ws.Cells.Item(15, 3)
$cell.NumberFormat = "@"
$cell.Value = "Facility5"
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "C"
$ws.Cells.Item(15, 5).Value = 231
$ws.Cells.Item(15, 6).Value = 3003
$cell = $ws.Cells.Item(15, 7)
$cell.NumberFormat = "@"
$cell.Value = "B"
$ws.Cells.Item(15, 8).Value = 13
$cell = $ws.Cells.Item(15, 9)
$cell.NumberFormat = "@"
$cell.Value = "3%"
$ws.Cells.Item(15, 10).Value = 12.61
$ws.Cells.Item(15, 11).Value = 163.93
$ws.Cells.Item(15, 12).Value = 13
$ws.Cells.Item(15, 13).Value = 2839.07
$cell = $ws.Cells.Item(15, 14)
$cell.NumberFormat = "@"
$cell.Value = "5%"
$ws.Cells.Item(15, 15).Value = 8.1965

# Update the LP Model sheet: the custom rule Rule_0_1 changes from a simple
# capacity-style bound to a balance constraint across awarded suppliers.
$lp = $wb.Worksheets.Item("LP Model")
$cell = $lp.Cells.Item(2, 1)
$text = $cell.Value
$oldLine = "Rule_0_1: x_B_1 + x_C_1 >= 500"
$newLine = "Rule_0_1: - 0.5 x_A_1 + 0.5 x_B_1 - 0.5 x_C_1 >= 0"
$text = $text.Replace($oldLine, $newLine)
$cell.Value = $text
